# feat: add 2022-Q1 data
#
# 1. Create a new worksheet "2022-Q1" (fund-holding detail, same shape as the
#    other quarterly sheets) positioned right before "总计".
# 2. Insert a new leading row into "总计" summarising the 2022-Q1 quarter and
#    renumber the helper index column.

$wb = $excel.ActiveWorkbook

# Helper: write a value that LOOKS like a number (e.g. "010746", "1.16")
# into a cell while forcing it to stay a text value, the same way the
# original workbook stored these "numeric-looking" strings.
function Set-TextValue($range, [string]$text) {
    $range.ClearFormats()
    $range.Value = "'" + $text
}

# ---------------------------------------------------------------------
# Step 1: build the new "2022-Q1" sheet by cloning "2021-Q4" (identical
# column layout/styling) right before "总计", then replace its data.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$beforeTarget = $wb.Worksheets.Item("总计")
$template.Copy($beforeTarget)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template has 7 data rows (2021-Q4 holds 7 funds); 2022-Q1 only has 4,
# so drop the trailing rows 6-8.
$newSheet.Range("A6:H8").Delete(-4162)

# Row 2 - 010746 富安达长三角区域主题混合
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "010746"
$newSheet.Range("C2").Value = "富安达长三角区域主题混合"
Set-TextValue $newSheet.Range("D2") "1.16"
Set-TextValue $newSheet.Range("E2") "92.45"
Set-TextValue $newSheet.Range("F2") "5.83"
Set-TextValue $newSheet.Range("G2") "0.0676"
$newSheet.Range("H2").Value = 3

# Row 3 - 009789 富安达科技创新混合
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "009789"
$newSheet.Range("C3").Value = "富安达科技创新混合"
Set-TextValue $newSheet.Range("D3") "0.56"
Set-TextValue $newSheet.Range("E3") "94.07"
Set-TextValue $newSheet.Range("F3") "3.78"
Set-TextValue $newSheet.Range("G3") "0.0212"
$newSheet.Range("H3").Value = 10

# Row 4 - 004549 富安达消费主题灵活配置混合
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "004549"
$newSheet.Range("C4").Value = "富安达消费主题灵活配置混合"
Set-TextValue $newSheet.Range("D4") "0.28"
Set-TextValue $newSheet.Range("E4") "93.19"
Set-TextValue $newSheet.Range("F4") "4.99"
Set-TextValue $newSheet.Range("G4") "0.0140"
$newSheet.Range("H4").Value = 8

# Row 5 - 005443 国金量化多策略灵活配置混合
$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet.Range("B5") "005443"
$newSheet.Range("C5").Value = "国金量化多策略灵活配置混合"
Set-TextValue $newSheet.Range("D5") "0.51"
Set-TextValue $newSheet.Range("E5") "64.10"
Set-TextValue $newSheet.Range("F5") "0.64"
Set-TextValue $newSheet.Range("G5") "0.0033"
$newSheet.Range("H5").Value = 10

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q1 row at the top of "总计" and renumber the
# index column (A) for the rows that got shifted down.
# (re-fetch the sheet by name: inserting/copying sheets can shift what
# an old worksheet reference points at)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert(-4121)

# New row 2 content - no inherited border/bold formatting on B:D.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.11

# A2 should use the same "index" style as the rest of column A.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0

# Renumber A3:A7 (previously A2:A6 = 0..4) to 1..5.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

$excel.CutCopyMode = 0
